# Auto-generated edit script: updates crypto price/volume table (Sun Apr 7 17:36:40 UTC 2024 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to round-trip numeric-looking strings through a text formula
# result (Copy + PasteSpecial values) so the D-column keeps its original "General"
# (text) storage instead of Excel auto-converting "582.54" etc. into a real number.
$helper = $ws.Range("Z1")

$helper.Formula = "=`"69.824.06`""
$helper.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  +2.51%  "
$helper.Formula = "=`"3.381.17`""
$helper.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  +1.28%  "
$helper.Formula = "=`"582.54`""
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.33%  "
$helper.Formula = "=`"180.34`""
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("E7").Value = "  +0.06%  "
$helper.Formula = "=`"0.595`""
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  +9.00%  "
$helper.Formula = "=`"0.591`""
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +1.30%  "
$helper.Formula = "=`"48.62`""
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +1.25%  "
$helper.Formula = "=`"0.0000286`""
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +4.59%  "
$helper.Formula = "=`"682.54`""
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -1.86%  "
$helper.Formula = "=`"8.63`""
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +2.25%  "
$helper.Formula = "=`"3.930.85`""
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +1.27%  "
$helper.Formula = "=`"69.806.46`""
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("E17").Value = "  +0.95%  "
$helper.Formula = "=`"3.387.28`""
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -0.26%  "
$helper.Formula = "=`"17.72`""
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +1.27%  "
$helper.Formula = "=`"11.32`""
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +1.27%  "
$helper.Formula = "=`"0.914`""
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +1.98%  "
$helper.Formula = "=`"17.34`""
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("E23").Value = "  -1.28%  "
$helper.Formula = "=`"101.99`""
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +0.47%  "
$helper.Formula = "=`"9.86`""
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +3.84%  "
$helper.Formula = "=`"33.66`""
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +1.58%  "
$helper.Formula = "=`"8.78`""
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +2.84%  "
$helper.Formula = "=`"6.94`""
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -0.36%  "
$helper.Formula = "=`"3.85`""
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +15.21%  "
$helper.Formula = "=`"11.10`""
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +0.31%  "
$helper.Formula = "=`"557.81`""
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -1.79%  "
$helper.Formula = "=`"0.106`""
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +0.91%  "
$helper.Formula = "=`"58.09`""
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("E36").Value = "  +0.03%  "
$helper.Formula = "=`"3.614.27`""
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("E38").Value = "  +2.78%  "
$helper.Formula = "=`"35.43`""
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +0.79%  "
$helper.Formula = "=`"0.0₃0733`""
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +8.69%  "
$helper.Formula = "=`"2.79`""
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +6.39%  "
$helper.Formula = "=`"3.34`""
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +5.16%  "
$ws.Range("E43").Value = "  +3.92%  "
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$helper.Formula = "=`"1.00`""
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$helper.Formula = "=`"1.38`""
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +3.70%  "
$helper.Formula = "=`"130.62`""
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  +0.93%  "

$helper.Clear()
$excel.CutCopyMode = $false
